$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.956.12"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3
$ws.Range("D3").Value = "3.498.87"
$ws.Range("E3").Value = "  -1.26%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.62%  "

# Row 7
$ws.Range("E7").Value = "  -1.63%  "

# Row 8
$ws.Range("D8").Value = "3.493.16"
$ws.Range("E8").Value = "  -1.31%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.193"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.29%  "

# Row 11
$ws.Range("E11").Value = "  +6.74%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.585"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "45.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.10%  "

# Row 14
$ws.Range("E14").Value = "  -2.03%  "

# Row 15
$ws.Range("D15").Value = "4.065.25"
$ws.Range("E15").Value = "  -0.91%  "

# Row 16
$ws.Range("E16").Value = "  -0.90%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "612.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.35%  "

# Row 18
$ws.Range("D18").Value = "3.507.66"
$ws.Range("E18").Value = "  -0.91%  "

# Row 19
$ws.Range("D19").Value = "69.945.41"
$ws.Range("E19").Value = "  -0.93%  "

# Row 20
$ws.Range("E20").Value = "  +0.90%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.21%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.876"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.22%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.44%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "98.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.46%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "15.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.30%  "

# Row 26
$ws.Range("E26").Value = "  -3.75%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("E28").Value = "  -2.52%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.74%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.97"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.19%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.32%  "

# Row 33
$ws.Range("E33").Value = "  -4.83%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.36%  "

# Row 35
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "626.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.46%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0995"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.79%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0480"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.29%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "56.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.79%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.145"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.32%  "

# Row 43
$ws.Range("D43").Value = "3.349.31"
$ws.Range("E43").Value = "  +0.31%  "

# Row 44
$ws.Range("E44").Value = "  +2.02%  "

# Row 45
$ws.Range("E45").Value = "  -5.99%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.40%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "31.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.25%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.45%  "

# Row 49
$ws.Range("E49").Value = "  +0.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.66%  "
